# The "vip" column (F) was stored as the text string "0" for every data row
# (F2:F101). Re-write it as a true numeric 0 so the column is numeric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 101; $row++) {
    $ws.Cells.Item($row, 6).Value = 0
}

# A handful of rows had a blank/corrupt "duration" (E) cell (stored as an
# empty string cell with no value). Fill them in with the numeric value 57.
$durationRows = @(4, 9, 13, 14, 15, 36, 47, 49, 51, 52, 53, 56, 64, 67, 69, 72, 75, 78, 80, 84, 88, 90, 92, 97)
foreach ($row in $durationRows) {
    $ws.Cells.Item($row, 5).Value = 57
}
